$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (Price and Volume(1h) columns) per commit diff
$updates = @{
    "D2" = "275.01"
    "E2" = "-2.17%"
    "D3" = "27.14"
    "E3" = "1.11%"
    "D4" = "4.753"
    "E4" = "-3.81%"
    "D5" = "0.06302"
    "E5" = "-1.78%"
    "E6" = "-0.95%"
    "D7" = "1.377"
    "E7" = "36.71%"
    "D8" = "0.8763"
    "E8" = "-1.04%"
    "E9" = "1.63%"
    "D10" = "0.05070"
    "E10" = "-1.87%"
    "D11" = "0.07630"
    "E11" = "3.04%"
    "D12" = "0.02970"
    "E12" = "-4.26%"
    "D13" = "0.09002"
    "E13" = "-0.60%"
    "E14" = "0.09%"
    "D15" = "0.0006366"
    "E15" = "1.04%"
    "D16" = "0.005977"
    "E16" = "-0.98%"
    "D17" = "3.445"
    "E17" = "-1.80%"
    "D18" = "3.295"
    "E18" = "-1.74%"
    "E19" = "-1.16%"
    "E20" = "0.13%"
    "E21" = "-0.38%"
    "E22" = "-0.65%"
    "D23" = "0.04393"
    "E23" = "0.92%"
    "D24" = "0.001169"
    "E24" = "-0.50%"
    "D25" = "0.003840"
    "E25" = "4.18%"
    "D26" = "0.0001199"
    "E26" = "-0.01%"
    "D27" = "0.0001935"
    "E27" = "14.29%"
    "D40" = "0.04097"
    "E40" = "-0.09%"
    "D41" = "0.006817"
    "E41" = "2.15%"
    "D42" = "0.1172"
    "E42" = "-0.50%"
    "D43" = "0.002099"
    "E43" = "-11.03%"
    "E44" = "-11.79%"
    "D45" = "0.00005177"
    "E45" = "-1.53%"
    "D46" = "1.486"
    "E46" = "-36.89%"
    "E47" = "2.34%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
